$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New table data (rows 16-28), columns C..G
# Columns: C = N Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora, F = Valor Mora, G = Salario Basico
$data = @(
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2404", 52000, 1300000),
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2405", 52000, 1300000),
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2406", 52000, 1300000),
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2407", 52000, 1300000),
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2408", 52000, 1300000),
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2409", 52000, 1300000),
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2410", 52000, 1300000),
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2411", 52000, 1300000),
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2412", 52000, 1300000),
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2501", 52000, 1300000),
    @("1128050520", "RANDY JAVIER TORRENTE HANNA",     "2502", 32933, 1300000),
    @("33102541",   "MAIRA ALEJANDRA MARTINEZ ZUÑIGA", "2502", 32933, 1300000),
    @("1041977150", "JORGE HUMBERTO RAMIREZ MARTINEZ", "2502", 32933, 1300000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 3).Value = $rowData[0]   # C - N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $rowData[1]   # D - Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $rowData[2]   # E - Periodo Mora
    $ws.Cells.Item($row, 6).Value = $rowData[3]   # F - Valor Mora
    $ws.Cells.Item($row, 7).Value = $rowData[4]   # G - Salario Basico
}

$wb.Save()
